# Wed, May 06, 2020  6:06:29 AM
#
# 1) Slide 6's table (graphic frame "Google Shape;127;p18") gets a new
#    built-in table style applied (was the deck's sole custom style
#    {2E2B957D-92FA-425C-9346-58970B3CACA9}, becomes the built-in
#    "Medium Style 2 - Accent 1" {3BA586D5-DD53-4A56-97B9-7F8A089B7527}).
#
# 2) The slide master's theme colour scheme is reset from the custom
#    "Integral" palette back to the stock Office default palette.

$p = $ppt.ActivePresentation

# --- (1) table style -------------------------------------------------
$slide = $p.Slides.Item(6)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{3BA586D5-DD53-4A56-97B9-7F8A089B7527}")

# --- (2) theme colours -------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function Set-SchemeRGB([int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $colorScheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# dk1, lt1, dk2, lt2, accent1..accent6, hlink, folHlink
Set-SchemeRGB 1  "000000"
Set-SchemeRGB 2  "FFFFFF"
Set-SchemeRGB 3  "44546A"
Set-SchemeRGB 4  "E7E6E6"
Set-SchemeRGB 5  "5B9BD5"
Set-SchemeRGB 6  "ED7D31"
Set-SchemeRGB 7  "A5A5A5"
Set-SchemeRGB 8  "FFC000"
Set-SchemeRGB 9  "4472C4"
Set-SchemeRGB 10 "70AD47"
Set-SchemeRGB 11 "0563C1"
Set-SchemeRGB 12 "954F72"

Write-Host "Table style + theme colours updated"
